$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 11 - "Nhap hoi vien" task: mark progress 100% and note UI change
$ws.Range("F11").Value = 1
$ws.Range("F11").NumberFormat = "0%"
$ws.Range("G11").Value = "Có thay đổi giao diện"

# Row 12 - "Cap nhat thong tin hoi vien" task: mark progress 100% and note UI change
$ws.Range("F12").Value = 1
$ws.Range("F12").NumberFormat = "0%"
$ws.Range("G12").Value = "Có thay đổi giao diện"

# Copy G11's formatting down onto G12 (fill-handle style drag-down)
$ws.Range("G11").Copy()
$ws.Range("G12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Move the active selection to C11
$ws.Range("C11").Select()
